$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 121, pushing existing rows 121-130 down to 122-131
$ws.Rows("121:121").Insert()

# Populate the new row 121 with the weekly price-report entry
$ws.Cells.Item(121, 1).Value = 3
$ws.Cells.Item(121, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(121, 3).Value = "Coquimbo"
$ws.Cells.Item(121, 4).Value = 44585
$ws.Cells.Item(121, 5).Value = 5
$ws.Cells.Item(121, 6).Value = 100112030
$ws.Cells.Item(121, 7).Value = "Poroto granado"
$ws.Cells.Item(121, 8).Value = "Sin especificar"
$ws.Cells.Item(121, 9).Value = "Primera"
$ws.Cells.Item(121, 10).Value = 38
$ws.Cells.Item(121, 11).Value = 23000
$ws.Cells.Item(121, 12).Value = 23000
$ws.Cells.Item(121, 13).Value = 23000
$ws.Cells.Item(121, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(121, 15).Value = "Provincia de Petorca"
$ws.Cells.Item(121, 16).Value = 920
$ws.Cells.Item(121, 17).Value = 25
$ws.Cells.Item(121, 18).Value = "Hortaliza"
